$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first column (A) held a running index (1, 17) that is no longer
# wanted. Deleting the entire column shifts B:F left into A:E, matching
# the target layout (headers now start at A1, data columns follow).
$ws.Range("A1").EntireColumn.Delete()
